$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: force-merge a paragraph's runs into a single run holding $newText.
# A plain one-shot assignment is a no-op when the replacement text is
# byte-identical to the paragraph's current (multi-run) text, so we first
# stomp the range with a placeholder and then write the real text — this
# guarantees Word actually rewrites the run structure (dropping proofErr
# markers and collapsing the split runs) while keeping the formatting of
# the paragraph's first original run.
# ---------------------------------------------------------------------------
function Merge-ParagraphText($idx, $newText) {
    $p = $d.Paragraphs($idx)
    $start = $p.Range.Start
    $end = $p.Range.End - 1   # exclude the trailing paragraph mark
    $r = $d.Range($start, $end)
    $r.Text = "@@TMP@@"
    $r2 = $d.Range($start, $start + 7)
    $r2.Text = $newText
}

# 1) March 11 paragraph — collapse the gramStart/gramEnd-split runs.
Merge-ParagraphText 1 "March 11, 2020 - WHO declares the novel coronavirus outbreak to be a pandemic. WHO says the outbreak is the first pandemic caused by a coronavirus. In an Oval Office address, Trump announces that he is restricting travel from Europe to the United States for 30 days in an attempt to slow the spread of coronavirus. The ban, which applies to the 26 countries in the Schengen Area, applies only to foreign nationals and not American citizens and permanent residents who'd be screened before entering the"

# 2) The May 11 / April 8 / April 15 items (currently paragraphs 5, 6, 7)
#    each get their split runs collapsed into one, in place, before being
#    reordered below.
Merge-ParagraphText 5 "May 11, 2020 - Trump and his administration announce that the federal government is sending `$11 billion to states to expand coronavirus testing capabilities. The relief package signed on April 24 includes `$25 billion for testing, with `$11 billion for states, localities, territories and tribes"

Merge-ParagraphText 6 "April 8 — Troubles With the COVID-19 Cocktail “What do you have to lose?” Trump asks when touting the malaria drug hydroxychloroquine or the related chloroquine as possible treatments for COVID-19. With a common antibiotic, azithromycin, the drug cocktail becomes an early candidate to prevent hospitalization or death. But Trump’s promotion of the combination, despite known heart risks for some patients, prompts the American Heart Association, the American College of Cardiology, and the Heart Rhythm Society to warn in a joint guidance that the drugs are not for everyone."

Merge-ParagraphText 7 "April 15 - Governor Andrew M. Cuomo today announced he will issue an Executive Order requiring all people in New York to wear a mask or a face covering when out in public and in situations where social distancing cannot be maintained, such as on public transportation. The Executive Order will go into effect on Friday, April 17th."

# 3) Reorder paragraphs 5-7 from (May 11, April 8, April 15) to
#    (April 8, April 15, May 11). Each paragraph now holds exactly one run,
#    so a FormattedText copy carries the right run properties along.
#    A cyclic rotation needs a scratch slot because FormattedText handles
#    stay bound to their source range; we physically park the May 11
#    content in a temporary trailing paragraph while the other two shift up.
$endIdx = $d.Content.End
$scratchRange = $d.Range($endIdx - 1, $endIdx - 1)
$scratchRange.InsertParagraphAfter()
$scratchIdx = $d.Paragraphs.Count

$p5 = $d.Paragraphs(5)
$scratch = $d.Paragraphs($scratchIdx)
$scratch.Range.FormattedText = $p5.Range.FormattedText   # scratch <- May 11

$p6 = $d.Paragraphs(6)
$d.Paragraphs(5).Range.FormattedText = $p6.Range.FormattedText   # P5 <- April 8

$p7 = $d.Paragraphs(7)
$d.Paragraphs(6).Range.FormattedText = $p7.Range.FormattedText   # P6 <- April 15

$scratch2 = $d.Paragraphs($scratchIdx)
$d.Paragraphs(7).Range.FormattedText = $scratch2.Range.FormattedText   # P7 <- May 11 (from scratch)

$scratch3 = $d.Paragraphs($scratchIdx)
$delRange = $d.Range($scratch3.Range.Start, $scratch3.Range.End)
$delRange.Delete()

# 4) August 23, 2021 paragraph (now still paragraph 12 — unaffected by the
#    reorder above, which only touched paragraphs 5-7) — collapse gramStart
#    split run.
Merge-ParagraphText 12 "August 23, 2021 - The US FDA grants full approval to the Pfizer/BioNTech Covid-19 vaccine for people age 16 and older, making it the first coronavirus vaccine approved by the FDA"

# 5) November 2, 2021 paragraph (still paragraph 13) — collapse spellStart
#    split run around "Walensky".
Merge-ParagraphText 13 "November 2, 2021 - Walensky says she is endorsing a recommendation to vaccinate children ages 5-11 against Covid-19, clearing the way for immediate vaccination of the youngest age group yet in the US"

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    Write-Output "$i : $($t.Substring(0, [Math]::Min(70, $t.Length)))"
}
